$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO"
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M10").Value = 7046.97
$ws1.Range("K12").Value = 1373.9
$ws1.Range("M12").Value = 2522.28
$ws1.Range("K18").Value = "1 de 16"

# Sheet "VENTA MENSUAL"
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F10").Value = 7662.57
$ws2.Range("F12").Value = 3896.18
$ws2.Range("F18").Value = 17222.84
